# Updated symbol list on Thu Dec 15 08:46:27 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as text (they were written as plain
# numeric-looking strings, not real numbers), so each new value is entered
# with a leading apostrophe to force text entry and keep exact formatting
# (this also preserves trailing zeros, e.g. row 19: 0.006210).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'265.28"
$ws.Range("D4").Value  = "'6.205"
$ws.Range("D5").Value  = "'0.06146"
$ws.Range("D6").Value  = "'3.564"
$ws.Range("D7").Value  = "'6.712"
$ws.Range("D8").Value  = "'1.348"
$ws.Range("D9").Value  = "'0.8268"
$ws.Range("D10").Value = "'0.01357"
$ws.Range("D12").Value = "'0.08217"
$ws.Range("D13").Value = "'0.03397"
$ws.Range("D14").Value = "'0.03151"
$ws.Range("D15").Value = "'0.09237"
$ws.Range("D16").Value = "'3.891"
$ws.Range("D17").Value = "'0.001694"
$ws.Range("D18").Value = "'0.04794"
$ws.Range("D19").Value = "'0.006210"
$ws.Range("D20").Value = "'0.006305"
$ws.Range("D21").Value = "'0.001098"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D23").Value = "'3.744"
$ws.Range("D24").Value = "'2.301"
$ws.Range("D27").Value = "'0.0002684"
$ws.Range("D40").Value = "'0.04612"
$ws.Range("D41").Value = "'0.006966"
$ws.Range("D42").Value = "'0.1135"
$ws.Range("D43").Value = "'0.003134"
$ws.Range("D44").Value = "'0.01088"
$ws.Range("D45").Value = "'0.00006169"

$ws.Range("D47").Value = "'0.7707"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").Value = "'0.2047"

$ws.Range("D49").Value = "'0.00001501"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"
